$d = $word.ActiveDocument

# Header(2) == "first page" header (header1.xml) -> BTec logo
#   name: image1.jpg -> image2.jpg
$h2 = $d.Sections(1).Headers(2)
if ($h2.Exists -and $h2.Range.InlineShapes.Count -ge 1) {
    $h2.Range.InlineShapes(1).Name = "image2.jpg"
}

# Footer(1) == default footer (footer2.xml) -> Pearson logo (docPr id=2)
#   name: image2.png -> image1.png
$f1 = $d.Sections(1).Footers(1)
if ($f1.Exists -and $f1.Range.InlineShapes.Count -ge 1) {
    $f1.Range.InlineShapes(1).Name = "image1.png"
}

# Footer(2) == "first page" footer (footer1.xml) -> Pearson logo (docPr id=3)
#   name: image2.png -> image1.png
$f2 = $d.Sections(1).Footers(2)
if ($f2.Exists -and $f2.Range.InlineShapes.Count -ge 1) {
    $f2.Range.InlineShapes(1).Name = "image1.png"
}
